$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Intermediate Data Visualization with Seaborn"
$ws.Range("J2").Value = 4
